# Update "想去人数" (F column) figures across sheets following the
# regenerated output referenced in the commit message.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 140
$ws1.Range("F7").Value = 301
$ws1.Range("F9").Value = 2067
$ws1.Range("F10").Value = 361
$ws1.Range("F11").Value = 4977
$ws1.Range("F12").Value = 101
$ws1.Range("F13").Value = 347

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 9
$ws2.Range("F5").Value = 14

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 140
$ws4.Range("F9").Value = 301
$ws4.Range("F11").Value = 9
$ws4.Range("F12").Value = 14
$ws4.Range("F13").Value = 2067
$ws4.Range("F14").Value = 361
$ws4.Range("F15").Value = 4977
$ws4.Range("F16").Value = 101
$ws4.Range("F17").Value = 347
